{"js": "// Ordered (old, new) text pairs taken from the canonical-OOXML diff.\n// Every old value is unique in the document, so a literal search+replace\n// for each pair (in any order) reproduces the target edit exactly.\nconst pairs = [\n  [\"2024-11-29 Friday\", \"2024-11-30 Saturday\"],\n  [\"89-5=84\", \"90+9=99\"],\n  [\"77-73=4\", \"25+29=54\"],\n  [\"78-57=21\", \"20+45=65\"],\n  [\"89-34=55\", \"0+41=41\"],\n  [\"61+22=83\", \"5+8=13\"],\n  [\"2+94=96\", \"96-84=12\"],\n  [\"88-73=15\", \"42+37=79\"],\n  [\"18+81=99\", \"68+10=78\"],\n  [\"88-19=69\", \"59-21=38\"],\n  [\"96-82=14\", \"36+37=73\"],\n  [\"87-74=13\", \"37+6=43\"],\n  [\"92-69=23\", \"50+24=74\"],\n  [\"76-38=38\", \"30+43=73\"],\n  [\"58+20=78\", \"30+9=39\"],\n  [\"17+39=56\", \"7+59=66\"],\n  [\"73-30=43\", \"95-65=30\"],\n  [\"76-5=71\", \"93-85=8\"],\n  [\"43+35=78\", \"54-24=30\"],\n  [\"24+39=63\", \"40+37=77\"],\n  [\"32+50=82\", \"61+13=74\"],\n  [\"22-12=10\", \"43-11=32\"],\n  [\"70+28=98\", \"15+8=23\"],\n  [\"39+11=50\", \"20+19=39\"],\n  [\"69-23=46\", \"79-21=58\"],\n  [\"43+56=99\", \"17-2=15\"],\n  [\"82+7=89\", \"75-5=70\"],\n  [\"49+31=80\", \"63-41=22\"],\n  [\"13+27=40\", \"94-29=65\"],\n  [\"78-43=35\", \"61+36=97\"],\n  [\"7+45=52\", \"21-4=17\"],\n  [\"59-56=3\", \"53+2=55\"],\n  [\"81-29=52\", \"79-0=79\"],\n  [\"84-57=27\", \"11+50=61\"],\n  [\"54+9=63\", \"56-19=37\"],\n  [\"58+23=81\", \"95-88=7\"],\n  [\"73-63=10\", \"81-35=46\"],\n  [\"65-40=25\", \"58-26=32\"],\n  [\"36+18=54\", \"50+33=83\"],\n  [\"78-59=19\", \"1+92=93\"],\n  [\"62-13=49\", \"85+0=85\"],\n  [\"16+42=58\", \"11+40=51\"],\n  [\"67-52=15\", \"55+34=89\"],\n  [\"70+25=95\", \"43+37=80\"],\n  [\"71-63=8\", \"17+9=26\"],\n  [\"18+0=18\", \"10+8=18\"],\n  [\"46-31=15\", \"20+13=33\"],\n  [\"38+42=80\", \"20+40=60\"],\n  [\"52+43=95\", \"51-13=38\"],\n  [\"29-0=29\", \"60-30=30\"],\n  [\"89-58=31\", \"13-0=13\"],\n  [\"56-0=56\", \"70-7=63\"],\n  [\"98-31=67\", \"38+33=71\"],\n  [\"42+28=70\", \"69-13=56\"],\n  [\"10+23=33\", \"90-61=29\"],\n  [\"57-32=25\", \"36+21=57\"],\n  [\"84-56=28\", \"7+41=48\"],\n  [\"24+59=83\", \"3+93=96\"],\n  [\"60+39=99\", \"90-41=49\"],\n  [\"90-19=71\", \"75+6=81\"],\n  [\"40+28=68\", \"2+80=82\"],\n  [\"49-27=22\", \"76-75=1\"],\n  [\"62-52=10\", \"76-59=17\"],\n  [\"17+37=54\", \"80-45=35\"],\n  [\"81-34=47\", \"13+18=31\"],\n  [\"24+50=74\", \"42+15=57\"],\n  [\"14+41=55\", \"53-14=39\"],\n  [\"85-82=3\", \"92-86=6\"],\n  [\"10+67=77\", \"97-61=36\"],\n  [\"22-16=6\", \"40+35=75\"],\n  [\"24-22=2\", \"28+32=60\"],\n  [\"86-47=39\", \"47+44=91\"],\n  [\"29+31=60\", \"46-37=9\"],\n  [\"60-50=10\", \"62-18=44\"],\n  [\"44-10=34\", \"40-39=1\"],\n  [\"59-42=17\", \"31+42=73\"],\n  [\"33+28=61\", \"65-42=23\"],\n  [\"63-27=36\", \"5+79=84\"],\n  [\"22+16=38\", \"11+18=29\"],\n  [\"87-38=49\", \"30+45=75\"],\n  [\"86-51=35\", \"47+49=96\"],\n  [\"86-37=49\", \"68-44=24\"],\n  [\"77-19=58\", \"46+32=78\"],\n  [\"81+8=89\", \"65-15=50\"],\n  [\"77-42=35\", \"73-53=20\"],\n  [\"86-32=54\", \"38-14=24\"],\n  [\"14-0=14\", \"23+46=69\"],\n  [\"28-6=22\", \"21+24=45\"],\n  [\"64-42=22\", \"26+52=78\"],\n  [\"97-89=8\", \"19+63=82\"],\n  [\"96-37=59\", \"31+43=74\"],\n  [\"88-2=86\", \"50+48=98\"],\n  [\"0+24=24\", \"69-24=45\"],\n  [\"19+62=81\", \"30-1=29\"],\n  [\"62-15=47\", \"17+22=39\"],\n  [\"13+85=98\", \"49+18=67\"],\n  [\"74-24=50\", \"7+15=22\"],\n  [\"13+84=97\", \"23+61=84\"],\n  [\"42-8=34\", \"44-3=41\"],\n  [\"63-45=18\", \"75+21=96\"],\n  [\"32-11=21\", \"62-20=42\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-11-29 Friday\", \"2024-11-30 Saturday\"),\n  @(\"89-5=84\", \"90+9=99\"),\n  @(\"77-73=4\", \"25+29=54\"),\n  @(\"78-57=21\", \"20+45=65\"),\n  @(\"89-34=55\", \"0+41=41\"),\n  @(\"61+22=83\", \"5+8=13\"),\n  @(\"2+94=96\", \"96-84=12\"),\n  @(\"88-73=15\", \"42+37=79\"),\n  @(\"18+81=99\", \"68+10=78\"),\n  @(\"88-19=69\", \"59-21=38\"),\n  @(\"96-82=14\", \"36+37=73\"),\n  @(\"87-74=13\", \"37+6=43\"),\n  @(\"92-69=23\", \"50+24=74\"),\n  @(\"76-38=38\", \"30+43=73\"),\n  @(\"58+20=78\", \"30+9=39\"),\n  @(\"17+39=56\", \"7+59=66\"),\n  @(\"73-30=43\", \"95-65=30\"),\n  @(\"76-5=71\", \"93-85=8\"),\n  @(\"43+35=78\", \"54-24=30\"),\n  @(\"24+39=63\", \"40+37=77\"),\n  @(\"32+50=82\", \"61+13=74\"),\n  @(\"22-12=10\", \"43-11=32\"),\n  @(\"70+28=98\", \"15+8=23\"),\n  @(\"39+11=50\", \"20+19=39\"),\n  @(\"69-23=46\", \"79-21=58\"),\n  @(\"43+56=99\", \"17-2=15\"),\n  @(\"82+7=89\", \"75-5=70\"),\n  @(\"49+31=80\", \"63-41=22\"),\n  @(\"13+27=40\", \"94-29=65\"),\n  @(\"78-43=35\", \"61+36=97\"),\n  @(\"7+45=52\", \"21-4=17\"),\n  @(\"59-56=3\", \"53+2=55\"),\n  @(\"81-29=52\", \"79-0=79\"),\n  @(\"84-57=27\", \"11+50=61\"),\n  @(\"54+9=63\", \"56-19=37\"),\n  @(\"58+23=81\", \"95-88=7\"),\n  @(\"73-63=10\", \"81-35=46\"),\n  @(\"65-40=25\", \"58-26=32\"),\n  @(\"36+18=54\", \"50+33=83\"),\n  @(\"78-59=19\", \"1+92=93\"),\n  @(\"62-13=49\", \"85+0=85\"),\n  @(\"16+42=58\", \"11+40=51\"),\n  @(\"67-52=15\", \"55+34=89\"),\n  @(\"70+25=95\", \"43+37=80\"),\n  @(\"71-63=8\", \"17+9=26\"),\n  @(\"18+0=18\", \"10+8=18\"),\n  @(\"46-31=15\", \"20+13=33\"),\n  @(\"38+42=80\", \"20+40=60\"),\n  @(\"52+43=95\", \"51-13=38\"),\n  @(\"29-0=29\", \"60-30=30\"),\n  @(\"89-58=31\", \"13-0=13\"),\n  @(\"56-0=56\", \"70-7=63\"),\n  @(\"98-31=67\", \"38+33=71\"),\n  @(\"42+28=70\", \"69-13=56\"),\n  @(\"10+23=33\", \"90-61=29\"),\n  @(\"57-32=25\", \"36+21=57\"),\n  @(\"84-56=28\", \"7+41=48\"),\n  @(\"24+59=83\", \"3+93=96\"),\n  @(\"60+39=99\", \"90-41=49\"),\n  @(\"90-19=71\", \"75+6=81\"),\n  @(\"40+28=68\", \"2+80=82\"),\n  @(\"49-27=22\", \"76-75=1\"),\n  @(\"62-52=10\", \"76-59=17\"),\n  @(\"17+37=54\", \"80-45=35\"),\n  @(\"81-34=47\", \"13+18=31\"),\n  @(\"24+50=74\", \"42+15=57\"),\n  @(\"14+41=55\", \"53-14=39\"),\n  @(\"85-82=3\", \"92-86=6\"),\n  @(\"10+67=77\", \"97-61=36\"),\n  @(\"22-16=6\", \"40+35=75\"),\n  @(\"24-22=2\", \"28+32=60\"),\n  @(\"86-47=39\", \"47+44=91\"),\n  @(\"29+31=60\", \"46-37=9\"),\n  @(\"60-50=10\", \"62-18=44\"),\n  @(\"44-10=34\", \"40-39=1\"),\n  @(\"59-42=17\", \"31+42=73\"),\n  @(\"33+28=61\", \"65-42=23\"),\n  @(\"63-27=36\", \"5+79=84\"),\n  @(\"22+16=38\", \"11+18=29\"),\n  @(\"87-38=49\", \"30+45=75\"),\n  @(\"86-51=35\", \"47+49=96\"),\n  @(\"86-37=49\", \"68-44=24\"),\n  @(\"77-19=58\", \"46+32=78\"),\n  @(\"81+8=89\", \"65-15=50\"),\n  @(\"77-42=35\", \"73-53=20\"),\n  @(\"86-32=54\", \"38-14=24\"),\n  @(\"14-0=14\", \"23+46=69\"),\n  @(\"28-6=22\", \"21+24=45\"),\n  @(\"64-42=22\", \"26+52=78\"),\n  @(\"97-89=8\", \"19+63=82\"),\n  @(\"96-37=59\", \"31+43=74\"),\n  @(\"88-2=86\", \"50+48=98\"),\n  @(\"0+24=24\", \"69-24=45\"),\n  @(\"19+62=81\", \"30-1=29\"),\n  @(\"62-15=47\", \"17+22=39\"),\n  @(\"13+85=98\", \"49+18=67\"),\n  @(\"74-24=50\", \"7+15=22\"),\n  @(\"13+84=97\", \"23+61=84\"),\n  @(\"42-8=34\", \"44-3=41\"),\n  @(\"63-45=18\", \"75+21=96\"),\n  @(\"32-11=21\", \"62-20=42\")\n)\n\nforeach ($pair in $pairs) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $pair[0]\n  $find.Replacement.Text = $pair[1]\n  # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n  #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n  # MatchCase:=True and MatchWholeWord:=False so each unique expression is matched\n  # literally (the strings include \"+ - =\" which are not word characters).\n  $find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
